$wb = $excel.ActiveWorkbook

# --- Rename the existing (only) sheet from "Sheet3" to "Hotels" ---
$hotels = $wb.Worksheets.Item(1)
$hotels.Name = "Hotels"

# --- Add a new "Reservations" sheet right after "Hotels" ---
$res = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $hotels)
$res.Name = "Reservations"

# --- Headers ---
$res.Range("A1").Value = "Reservation Code"
$res.Range("B1").Value = "Hotel Code"
$res.Range("C1").Value = "Customer Name"
$res.Range("D1").Value = "Check-In Date"
$res.Range("E1").Value = "Check-Out Date"
$res.Range("F1").Value = "Number of Rooms"

# --- Data row (references Desert Mirage's hotel code from the Hotels sheet) ---
$hotelCode = $hotels.Range("A7").Value2

$res.Range("A2").Value = "J02"
$res.Range("B2").Value = $hotelCode
$res.Range("C2").Value = "Adatum Corporation"
$res.Range("D2").Value = "8/6/2024"
$res.Range("E2").Value = "11/6/2024"
$res.Range("F2").Value = 2

# --- Date formatting for check-in/out columns ---
$res.Range("D2:E2").NumberFormat = "m/d/yyyy"

# --- Column widths to fit content, like the original sheet ---
$res.Columns.Item(1).AutoFit() | Out-Null
$res.Columns.Item(2).AutoFit() | Out-Null
$res.Columns.Item(3).AutoFit() | Out-Null
$res.Columns.Item(4).AutoFit() | Out-Null
$res.Columns.Item(5).AutoFit() | Out-Null
$res.Columns.Item(6).AutoFit() | Out-Null

# --- Match the selection left by the original author: cursor on B2 of the
#     Reservations sheet, but the Hotels sheet remains the active tab ---
$res.Range("B2").Select()
$hotels.Activate()
